$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/1/2025  Through  9/7/2025"

# --- Cells changing from numeric to the text placeholder "0" (shared style s=13) ---
# Donor cell C14 already has style s=13 / shared-string "0"; Copy() clones style+type,
# reusing the existing style index instead of fabricating a new one.
$ws.Range("C14").Copy($ws.Range("F15"))
$ws.Range("C14").Copy($ws.Range("C18"))
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("C14").Copy($ws.Range("F27"))
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("C14").Copy($ws.Range("D33"))

# --- Cells changing from numeric to the text placeholder "***.*" (shared style s=13) ---
# Donor cell E14 already has style s=13 / shared-string "***.*"
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("E14").Copy($ws.Range("E33"))

# --- Cells changing from the text placeholder "0" to numeric (style s=14) ---
# Donor cell J14 already has style s=14 (plain integer count format)
$ws.Range("J14").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 1
$ws.Range("J14").Copy($ws.Range("C17"))
$ws.Range("C17").Value = 2
$ws.Range("J14").Copy($ws.Range("C23"))
$ws.Range("C23").Value = 1

# --- Plain value updates (style/type unchanged) ---
$ws.Range("F14").Value = 1
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -100
$ws.Range("J15").Value = 7
$ws.Range("K15").Value = -28.571428571428
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = -40
$ws.Range("I16").Value = 31
$ws.Range("J16").Value = 45
$ws.Range("K16").Value = -31.111111111111
$ws.Range("L16").Value = 24
$ws.Range("M16").Value = -29.545454545454
$ws.Range("N16").Value = -83.243243243243
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 78
$ws.Range("J17").Value = 113
$ws.Range("K17").Value = -30.973451327433
$ws.Range("L17").Value = -9.302325581395
$ws.Range("M17").Value = 21.875
$ws.Range("N17").Value = -49.677419354838
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 66.666666666666
$ws.Range("M18").Value = -9.523809523809
$ws.Range("N18").Value = -87.417218543046
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -66.666666666666
$ws.Range("F19").Value = 10
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = -41.176470588235
$ws.Range("I19").Value = 86
$ws.Range("J19").Value = 112
$ws.Range("K19").Value = -23.214285714285
$ws.Range("L19").Value = -34.351145038167
$ws.Range("M19").Value = 6.172839506172
$ws.Range("N19").Value = -40.277777777777
$ws.Range("D20").Value = 2
$ws.Range("J20").Value = 44
$ws.Range("K20").Value = -61.363636363636
$ws.Range("L20").Value = -37.037037037037
$ws.Range("N20").Value = -93.089430894308
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = -68.75
$ws.Range("F21").Value = 28
$ws.Range("G21").Value = 47
$ws.Range("H21").Value = -40.425531914893
$ws.Range("I21").Value = 257
$ws.Range("J21").Value = 359
$ws.Range("K21").Value = -28.412256267409
$ws.Range("L21").Value = -14.617940199335
$ws.Range("M21").Value = 1.181102362204
$ws.Range("N21").Value = -75.430210325047
$ws.Range("F23").Value = 2
$ws.Range("I23").Value = 16
$ws.Range("K23").Value = -15.789473684210
$ws.Range("L23").Value = 100
$ws.Range("M23").Value = -11.111111111111
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = -30.769230769230
$ws.Range("F24").Value = 20
$ws.Range("G24").Value = 57
$ws.Range("H24").Value = -64.912280701754
$ws.Range("I24").Value = 327
$ws.Range("J24").Value = 338
$ws.Range("K24").Value = -3.254437869822
$ws.Range("L24").Value = -2.095808383233
$ws.Range("M24").Value = 53.521126760563
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -28.571428571428
$ws.Range("F25").Value = 13
$ws.Range("H25").Value = -53.571428571428
$ws.Range("I25").Value = 168
$ws.Range("J25").Value = 166
$ws.Range("K25").Value = 1.204819277108
$ws.Range("L25").Value = -5.084745762711
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -16.666666666666
$ws.Range("F26").Value = 22
$ws.Range("G26").Value = 22
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 166
$ws.Range("J26").Value = 188
$ws.Range("K26").Value = -11.702127659574
$ws.Range("L26").Value = 2.469135802469
$ws.Range("M26").Value = -33.064516129032
$ws.Range("D27").Value = 2
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -100
$ws.Range("J27").Value = 10
$ws.Range("K27").Value = -40
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = -33.333333333333
$ws.Range("F29").Value = 1
$ws.Range("F30").Value = 1
